$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete row 4 (Neutrophils row) entirely - shifts nothing up since it's the last row
$ws.Rows.Item(4).Delete()

# Update row 2 (ECs -> MuSCs via Tgfa/Erbb4) with new TPM-derived values
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.5
$ws.Cells.Item(2, 7).Value = 0.4368315
$ws.Cells.Item(2, 8).Value = 0.873663
$ws.Cells.Item(2, 9).Value = 0.5629387977071691
$ws.Cells.Item(2, 10).Value = 0.5629387977071691
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.0108025
$ws.Cells.Item(2, 14).Value = 0.021605
$ws.Cells.Item(2, 15).Value = 1
$ws.Cells.Item(2, 16).Value = 1
$ws.Cells.Item(2, 17).Value = 0.004718872278749999
$ws.Cells.Item(2, 18).Value = 0.018875489115
$ws.Cells.Item(2, 19).Value = 0.5629387977071691
$ws.Cells.Item(2, 20).Value = 0.5629387977071691

# Update row 3 (MuSCs -> MuSCs via Tgfa/Erbb4) with new TPM-derived values
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.3391525
$ws.Cells.Item(3, 8).Value = 0.6783049999999999
$ws.Cells.Item(3, 9).Value = 0.4370612022928307
$ws.Cells.Item(3, 10).Value = 0.4370612022928307
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.0108025
$ws.Cells.Item(3, 14).Value = 0.021605
$ws.Cells.Item(3, 15).Value = 1
$ws.Cells.Item(3, 16).Value = 1
$ws.Cells.Item(3, 17).Value = 0.00366369488125
$ws.Cells.Item(3, 18).Value = 0.014654779525
$ws.Cells.Item(3, 19).Value = 0.4370612022928307
$ws.Cells.Item(3, 20).Value = 0.4370612022928307
